# edit.ps1
# Applies the "Updated date and version" commit to the AB snapshot document:
#   1. Bumps the title-page version number  2.1.RC3  ->  3.0.RC1
#   2. Bumps the title-page date            2014-03-13 -> 2014-04-11
#      (the historical entries inside the revision-history table are left
#       untouched, exactly as in the source diff)
#   3. Re-creates the "_GoBack" bookmark spanning the whole document body,
#      which is the natural side-effect of Word re-saving a document after
#      an edit (it also renumbers all following bookmark ids by +1, which
#      matches the diff exactly).

$d = $word.ActiveDocument

# --- 1. Title-page version number --------------------------------------
# wdReplaceOne (1) so only the first occurrence (title page) is touched;
# the same text appears again later inside the revision-history table and
# must stay "2.1.RC3" there.
$d.Content.Find.Execute(
    "2.1.RC3", $true, $false, $false, $false, $false,
    $true, 1, $false, "3.0.RC1", 1) | Out-Null

# --- 2. Title-page date ---------------------------------------------------
# Same idea: only replace the first (title page) occurrence of the date,
# leave the historical revision-table row as-is.
$d.Content.Find.Execute(
    "2014-03-13", $true, $false, $false, $false, $false,
    $true, 1, $false, "2014-04-11", 1) | Out-Null

# --- 3. "_GoBack" bookmark -------------------------------------------------
# Word maintains a hidden "_GoBack" bookmark that marks the span covering
# the most recent edits. After the edits above it should wrap the entire
# document body (from the very start through to the end of the content).
$bodyStart = 0
$bodyEnd = $d.Content.End
$goBackRange = $d.Range($bodyStart, $bodyEnd)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
